$wb = $excel.ActiveWorkbook

# Starting state: UK, Belgium
$ukSheet = $wb.Worksheets.Item("UK")

# 1) Duplicate UK, place the copy right after UK -> rename to "Denmark"
#    Order so far: UK, Denmark, Belgium
$ukSheet.Copy([System.Reflection.Missing]::Value, $ukSheet) | Out-Null
$denmark = $wb.Worksheets.Item("UK (2)")
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2005"
$denmark.Activate() | Out-Null
$denmark.Range("A1:XFD1048576").Select() | Out-Null

# 2) Duplicate UK, place the copy right after Belgium -> rename to "Sweden"
#    Order so far: UK, Denmark, Belgium, Sweden
$belgium = $wb.Worksheets.Item("Belgium")
$ukSheet.Copy([System.Reflection.Missing]::Value, $belgium) | Out-Null
$sweden = $wb.Worksheets.Item("UK (2)")
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2023"
$sweden.Activate() | Out-Null
$sweden.Range("A1:XFD1048576").Select() | Out-Null

# 3) Duplicate UK, place the copy right after Sweden -> rename to "Norway"
#    Final order: UK, Denmark, Belgium, Sweden, Norway
$ukSheet.Copy([System.Reflection.Missing]::Value, $sweden) | Out-Null
$norway = $wb.Worksheets.Item("UK (2)")
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1920"

# Norway ends up as the active/selected tab with B2:B4 selected
$norway.Activate() | Out-Null
$norway.Range("B2:B4").Select() | Out-Null
